$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (labels shift from Pin_2_x / Pin_3_0-2 to Pin_3_x / Pin_4_0-2)
$ws.Range("B1").Value = "Pin_3_0"
$ws.Range("C1").Value = "Pin_3_2"
$ws.Range("D1").Value = "Pin_3_4"
$ws.Range("E1").Value = "Pin_3_6"
$ws.Range("F1").Value = "Pin_3_8"
$ws.Range("G1").Value = "Pin_4_0"
$ws.Range("H1").Value = "Pin_4_2"

# Update data rows with new computed values
$ws.Range("B2").Value = 27289.41365097705
$ws.Range("C2").Value = 28997.03757680443
$ws.Range("D2").Value = 30699.11894048838
$ws.Range("E2").Value = 32395.85801142112
$ws.Range("F2").Value = 34087.45887438462
$ws.Range("G2").Value = 35774.12905413556
$ws.Range("H2").Value = 37456.07915945948

$ws.Range("B3").Value = 27343.41092200421
$ws.Range("C3").Value = 29053.96309513143
$ws.Range("D3").Value = 30758.90060446428
$ws.Range("E3").Value = 32458.42349011225
$ws.Range("F3").Value = 34152.73569285308
$ws.Range("G3").Value = 35842.0446720012
$ws.Range("H3").Value = 37526.5610430948

$ws.Range("B4").Value = 27397.72999821576
$ws.Range("C4").Value = 29111.2251950953
$ws.Range("D4").Value = 30819.03288010401
$ws.Range("E4").Value = 32521.35286902505
$ws.Range("F4").Value = 34218.388964855
$ws.Range("G4").Value = 35910.34856884402
$ws.Range("H4").Value = 37597.44231171037

$ws.Range("B5").Value = 27452.37408877807
$ws.Range("C5").Value = 29168.82720664944
$ws.Range("D5").Value = 30879.51920801507
$ws.Range("E5").Value = 32584.64968948758
$ws.Range("F5").Value = 34284.42232272403
$ws.Range("G5").Value = 35979.04445852413
$ws.Range("H5").Value = 37668.72675146379
